$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set every cell in column C from row 2 to row 252 to the constant value 7573
$ws.Range("C2:C252").Value = 7573
